$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# --- zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$mdFile = "070da431-3186-4a66-a450-cfe76cedcac5.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ce01f7ef809b523de82b792740744580a54f699/e2e/070da431-3186-4a66-a450-cfe76cedcac5.md"
$zhCnXlf = "070da431-3186-4a66-a450-cfe76cedcac5.0c2cface3e3b9341f5f1ceb2fefce35a925286e1.zh-cn.xlf"
$deDeXlf = "070da431-3186-4a66-a450-cfe76cedcac5.0c2cface3e3b9341f5f1ceb2fefce35a925286e1.de-de.xlf"

$wsZhCn.Range("I2").Value = $mdFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFile)
$wsZhCn.Range("J2").Value = $zhCnXlf
$wsZhCn.Range("K2").Value = "2016-09-05 07:13:31"

$wsZhCn.Range("I3").Value = $mdFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl, "", "", $mdFile)
$wsZhCn.Range("J3").Value = $zhCnXlf
$wsZhCn.Range("K3").Value = "2016-09-05 07:13:31"

# --- de-de sheet: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsDeDe.Range("I2").Value = $mdFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFile)
$wsDeDe.Range("J2").Value = $deDeXlf
$wsDeDe.Range("K2").Value = "2016-09-05 07:13:38"

$wsDeDe.Range("I3").Value = $mdFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl, "", "", $mdFile)
$wsDeDe.Range("J3").Value = $deDeXlf
$wsDeDe.Range("K3").Value = "2016-09-05 07:13:38"

# --- Column width updates ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40
